$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: category switches to "Amino acid metabolism" (freq 4), percentage recalculated
$ws.Range("A4").Value = "Amino acid metabolism"
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 23.529411764705884
$ws.Range("C4").NumberFormat = "0.0"

# Row 5: category switches to "Metabolism of cofactors and vitamins", frequency bumps 3 -> 4
$ws.Range("A5").Value = "Metabolism of cofactors and vitamins"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 23.529411764705884
$ws.Range("C5").NumberFormat = "0.0"

# Rows 6-9: same frequency (2), percentage recalculated against new total (17)
$ws.Range("C6").Value = 11.764705882352942
$ws.Range("C6").NumberFormat = "0.0"

$ws.Range("C7").Value = 11.764705882352942
$ws.Range("C7").NumberFormat = "0.0"

$ws.Range("C8").Value = 11.764705882352942
$ws.Range("C8").NumberFormat = "0.0"

$ws.Range("C9").Value = 11.764705882352942
$ws.Range("C9").NumberFormat = "0.0"

# Row 10: frequency stays 1, percentage recalculated against new total (17)
$ws.Range("C10").Value = 5.882352941176471
$ws.Range("C10").NumberFormat = "0.0"

# Selection moves to C4:C10 (the updated percentage column) with C4 active
$ws.Range("C4:C10").Select()
